$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 3).Value = 0.1040162398185203
$ws.Cells.Item(2, 4).Value = 0.1563759751366547
$ws.Cells.Item(2, 5).Value = 0.06340302046016433
$ws.Cells.Item(2, 6).Value = 5.606357340710986
$ws.Cells.Item(2, 7).Value = 0.002597308631189735
$ws.Cells.Item(2, 9).Value = 3.939253232387301
$ws.Cells.Item(2, 10).Value = 0.07528521835414992
$ws.Cells.Item(2, 13).Value = 3.144960152179095
$ws.Cells.Item(2, 14).Value = 1.433860535331263

# Row 3
$ws.Cells.Item(3, 3).Value = 0.1034179608638652
$ws.Cells.Item(3, 4).Value = 0.1359955768662928
$ws.Cells.Item(3, 5).Value = 0.05741139501736825
$ws.Cells.Item(3, 6).Value = 5.463541721084795
$ws.Cells.Item(3, 7).Value = 0.002610357156766684
$ws.Cells.Item(3, 9).Value = 3.809128486789888
$ws.Cells.Item(3, 10).Value = 0.07507091591273607
$ws.Cells.Item(3, 13).Value = 2.919836584983898
$ws.Cells.Item(3, 14).Value = 1.359489484099612

# Row 4
$ws.Cells.Item(4, 3).Value = 0.1031480886666429
$ws.Cells.Item(4, 4).Value = 0.1235491552488384
$ws.Cells.Item(4, 5).Value = 0.05371803391608765
$ws.Cells.Item(4, 6).Value = 5.38116670960045
$ws.Cells.Item(4, 7).Value = 0.002618754821234998
$ws.Cells.Item(4, 9).Value = 3.733041349383683
$ws.Cells.Item(4, 10).Value = 0.0749618873699518
$ws.Cells.Item(4, 13).Value = 2.782910756946421
$ws.Cells.Item(4, 14).Value = 1.314299558801338

# Row 5
$ws.Cells.Item(5, 3).Value = 0.1030622387928446
$ws.Cells.Item(5, 4).Value = 0.1184921867767059
$ws.Cells.Item(5, 5).Value = 0.05220908582172612
$ws.Cells.Item(5, 6).Value = 5.348907783197404
$ws.Cells.Item(5, 7).Value = 0.002622274519818062
$ws.Cells.Item(5, 9).Value = 3.702971898097275
$ws.Cells.Item(5, 10).Value = 0.07492317334317278
$ws.Cells.Item(5, 13).Value = 2.727430839532161
$ws.Cells.Item(5, 14).Value = 1.296005556872615

# Row 6
$ws.Cells.Item(6, 3).Value = 0.1030494295262372
$ws.Cells.Item(6, 4).Value = 0.1176533366563604
$ws.Cells.Item(6, 5).Value = 0.05195828554291282
$ws.Cells.Item(6, 6).Value = 5.343629538100799
$ws.Cells.Item(6, 7).Value = 0.002622864873489694
$ws.Cells.Item(6, 9).Value = 3.69803484748968
$ws.Cells.Item(6, 10).Value = 0.07491709144460046
$ws.Cells.Item(6, 13).Value = 2.718237430210564
$ws.Cells.Item(6, 14).Value = 1.292975236356909

# Row 7
$ws.Cells.Item(7, 3).Value = 0.1031468336836383
$ws.Cells.Item(7, 4).Value = 0.1234808966867718
$ws.Cells.Item(7, 5).Value = 0.05369769970002025
$ws.Cells.Item(7, 6).Value = 5.380726386917217
$ws.Cells.Item(7, 7).Value = 0.002618801893214469
$ws.Cells.Item(7, 9).Value = 3.732632057500865
$ws.Cells.Item(7, 10).Value = 0.07496134205643656
$ws.Cells.Item(7, 13).Value = 2.782161257088944
$ws.Cells.Item(7, 14).Value = 1.31405234578483

# Row 8
$ws.Cells.Item(8, 3).Value = 0.1037894870724472
$ws.Cells.Item(8, 4).Value = 0.1493337245208011
$ws.Cells.Item(8, 5).Value = 0.06133996911027495
$ws.Cells.Item(8, 6).Value = 5.55599512699132
$ws.Cells.Item(8, 7).Value = 0.002601728029317909
$ws.Cells.Item(8, 9).Value = 3.893582370239926
$ws.Cells.Item(8, 10).Value = 0.07520667243097279
$ws.Cells.Item(8, 13).Value = 3.067062134008268
$ws.Cells.Item(8, 14).Value = 1.408120137765792

# Row 9
$ws.Cells.Item(9, 3).Value = 0.1058411932937702
$ws.Cells.Item(9, 4).Value = 0.200652475904775
$ws.Cells.Item(9, 5).Value = 0.07622305992686762
$ws.Cells.Item(9, 6).Value = 5.94312752309321
$ws.Cells.Item(9, 7).Value = 0.002571281448740456
$ws.Cells.Item(9, 9).Value = 4.240445675180183
$ws.Cells.Item(9, 10).Value = 0.07586481423160052
$ws.Cells.Item(9, 13).Value = 3.636512385282828
$ws.Cells.Item(9, 14).Value = 1.596274709011055

# Row 10
$ws.Cells.Item(10, 3).Value = 0.1078571353104962
$ws.Cells.Item(10, 4).Value = 0.2388637426395235
$ws.Cells.Item(10, 5).Value = 0.08711209890235239
$ws.Cells.Item(10, 6).Value = 6.255866308602549
$ws.Cells.Item(10, 7).Value = 0.002550725738860195
$ws.Cells.Item(10, 9).Value = 4.515798889556351
$ws.Cells.Item(10, 10).Value = 0.07645379713208911
$ws.Cells.Item(10, 13).Value = 4.062124323232496
$ws.Cells.Item(10, 14).Value = 1.736680449797461

# Row 11
$ws.Cells.Item(11, 3).Value = 0.1088905950145147
$ws.Cells.Item(11, 4).Value = 0.2563864584628277
$ws.Cells.Item(11, 5).Value = 0.09205994204373269
$ws.Cells.Item(11, 6).Value = 6.404705072715899
$ws.Cells.Item(11, 7).Value = 0.002541760059101651
$ws.Cells.Item(11, 9).Value = 4.64585249908842
$ws.Cells.Item(11, 10).Value = 0.07674410726816916
$ws.Cells.Item(11, 13).Value = 4.257476107315995
$ws.Cells.Item(11, 14).Value = 1.801011256919026

# Row 12
$ws.Cells.Item(12, 3).Value = 0.1092991808392298
$ws.Cells.Item(12, 4).Value = 0.2630444509061931
$ws.Cells.Item(12, 5).Value = 0.09393307646472238
$ws.Cells.Item(12, 6).Value = 6.46204668632015
$ws.Cells.Item(12, 7).Value = 0.002538419730409913
$ws.Cells.Item(12, 9).Value = 4.695817855553088
$ws.Cells.Item(12, 10).Value = 0.07685721156074976
$ws.Cells.Item(12, 13).Value = 4.331713590831072
$ws.Cells.Item(12, 14).Value = 1.82543632773934

# Row 13
$ws.Cells.Item(13, 3).Value = 0.1092104097952529
$ws.Cells.Item(13, 4).Value = 0.2616094930231156
$ws.Cells.Item(13, 5).Value = 0.0935296809849433
$ws.Cells.Item(13, 6).Value = 6.449653012126191
$ws.Cells.Item(13, 7).Value = 0.00253913670349839
$ws.Cells.Item(13, 9).Value = 4.685024589697548
$ws.Cells.Item(13, 10).Value = 0.07683271230224875
$ws.Cells.Item(13, 13).Value = 4.315713370549389
$ws.Cells.Item(13, 14).Value = 1.820173113934743

# Row 14
$ws.Cells.Item(14, 3).Value = 0.1089238612322561
$ws.Cells.Item(14, 4).Value = 0.2569337503415738
$ws.Cells.Item(14, 5).Value = 0.09221405429616425
$ws.Cells.Item(14, 6).Value = 6.409402766728533
$ws.Cells.Item(14, 7).Value = 0.002541484154027273
$ws.Cells.Item(14, 9).Value = 4.649948650532593
$ws.Cells.Item(14, 10).Value = 0.07675334916805099
$ws.Cells.Item(14, 13).Value = 4.263578348052619
$ws.Cells.Item(14, 14).Value = 1.803019441096126

# Row 15
$ws.Cells.Item(15, 3).Value = 0.1087506018784978
$ws.Cells.Item(15, 4).Value = 0.2540727272294419
$ws.Cells.Item(15, 5).Value = 0.09140813827539773
$ws.Cells.Item(15, 6).Value = 6.384876949166767
$ws.Cells.Item(15, 7).Value = 0.002542929150733682
$ws.Cells.Item(15, 9).Value = 4.628557817467566
$ws.Cells.Item(15, 10).Value = 0.07670514837393938
$ws.Cells.Item(15, 13).Value = 4.23167862777774
$ws.Cells.Item(15, 14).Value = 1.792520658006282

# Row 16
$ws.Cells.Item(16, 3).Value = 0.1077919807532339
$ws.Cells.Item(16, 4).Value = 0.237721612973786
$ws.Cells.Item(16, 5).Value = 0.08678865798827218
$ws.Cells.Item(16, 6).Value = 6.246274517817483
$ws.Cells.Item(16, 7).Value = 0.002551319366372429
$ws.Cells.Item(16, 9).Value = 4.507398375940568
$ws.Cells.Item(16, 10).Value = 0.07643527091642
$ws.Cells.Item(16, 13).Value = 4.049393622308571
$ws.Cells.Item(16, 14).Value = 1.732485380812335

# Row 17
$ws.Cells.Item(17, 3).Value = 0.1072340592569816
$ws.Cells.Item(17, 4).Value = 0.2277282966244059
$ws.Cells.Item(17, 5).Value = 0.08395352505558407
$ws.Cells.Item(17, 6).Value = 6.162953892550092
$ws.Cells.Item(17, 7).Value = 0.002556564721799951
$ws.Cells.Item(17, 9).Value = 4.434318095053129
$ws.Cells.Item(17, 10).Value = 0.07627540764292462
$ws.Cells.Item(17, 13).Value = 3.938021338983305
$ws.Cells.Item(17, 14).Value = 1.695772150451774

# Row 18
$ws.Cells.Item(18, 3).Value = 0.1069240824781588
$ws.Cells.Item(18, 4).Value = 0.2219934571928945
$ws.Cells.Item(18, 5).Value = 0.08232228215242543
$ws.Cells.Item(18, 6).Value = 6.115646549049416
$ws.Cells.Item(18, 7).Value = 0.002559618007154014
$ws.Cells.Item(18, 9).Value = 4.392733782808733
$ws.Cells.Item(18, 10).Value = 0.07618556781029895
$ws.Cells.Item(18, 13).Value = 3.874125856244007
$ws.Cells.Item(18, 14).Value = 1.674699020272158

# Row 19
$ws.Cells.Item(19, 3).Value = 0.1068209911244651
$ws.Cells.Item(19, 4).Value = 0.2200539103855306
$ws.Cells.Item(19, 5).Value = 0.08176986793256447
$ws.Cells.Item(19, 6).Value = 6.099733945331764
$ws.Cells.Item(19, 7).Value = 0.002560658049425804
$ws.Cells.Item(19, 9).Value = 4.378730435754477
$ws.Cells.Item(19, 10).Value = 0.07615551331267412
$ws.Cells.Item(19, 13).Value = 3.852519568866398
$ws.Cells.Item(19, 14).Value = 1.667571518578313

# Row 20
$ws.Cells.Item(20, 3).Value = 0.1072923167090067
$ws.Cells.Item(20, 4).Value = 0.2287907346057523
$ws.Cells.Item(20, 5).Value = 0.08425538450197223
$ws.Cells.Item(20, 6).Value = 6.171759466535548
$ws.Cells.Item(20, 7).Value = 0.002556002592721395
$ws.Cells.Item(20, 9).Value = 4.442050894923966
$ws.Cells.Item(20, 10).Value = 0.07629220740314935
$ws.Cells.Item(20, 13).Value = 3.949860163497078
$ws.Cells.Item(20, 14).Value = 1.699675862563225

# Row 21
$ws.Cells.Item(21, 3).Value = 0.1090075556089545
$ws.Cells.Item(21, 4).Value = 0.2583064988928356
$ws.Cells.Item(21, 5).Value = 0.09260049699880568
$ws.Cells.Item(21, 6).Value = 6.421198371725552
$ws.Cells.Item(21, 7).Value = 0.002540793169103967
$ws.Cells.Item(21, 9).Value = 4.660231623928865
$ws.Cells.Item(21, 10).Value = 0.07677657439141683
$ws.Cells.Item(21, 13).Value = 4.278884458756977
$ws.Cells.Item(21, 14).Value = 1.808056156124167

# Row 22
$ws.Cells.Item(22, 3).Value = 0.1102292796274469
$ws.Cells.Item(22, 4).Value = 0.2777291810027691
$ws.Cells.Item(22, 5).Value = 0.0980517099489262
$ws.Cells.Item(22, 6).Value = 6.589947963321379
$ws.Cells.Item(22, 7).Value = 0.002531171937367235
$ws.Cells.Item(22, 9).Value = 4.807017369086196
$ws.Cells.Item(22, 10).Value = 0.07711159170928994
$ws.Cells.Item(22, 13).Value = 4.495454015622784
$ws.Cells.Item(22, 14).Value = 1.879263535868319

# Row 23
$ws.Cells.Item(23, 3).Value = 0.1095678370363231
$ws.Cells.Item(23, 4).Value = 0.2673500109582676
$ws.Cells.Item(23, 5).Value = 0.09514244715224152
$ws.Cells.Item(23, 6).Value = 6.49934737552735
$ws.Cells.Item(23, 7).Value = 0.002536277986311444
$ws.Cells.Item(23, 9).Value = 4.728282123414175
$ws.Cells.Item(23, 10).Value = 0.07693111392337215
$ws.Cells.Item(23, 13).Value = 4.379722363866534
$ws.Cells.Item(23, 14).Value = 1.841225078038264

# Row 24
$ws.Cells.Item(24, 3).Value = 0.1072659449761062
$ws.Cells.Item(24, 4).Value = 0.2283103740179513
$ws.Cells.Item(24, 5).Value = 0.08411891787841341
$ws.Cells.Item(24, 6).Value = 6.167776616481717
$ws.Cells.Item(24, 7).Value = 0.002556256613760589
$ws.Cells.Item(24, 9).Value = 4.438553557212828
$ws.Cells.Item(24, 10).Value = 0.07628460578764873
$ws.Cells.Item(24, 13).Value = 3.944507414680089
$ws.Cells.Item(24, 14).Value = 1.697910888924355

# Row 25
$ws.Cells.Item(25, 3).Value = 0.105198784713096
$ws.Cells.Item(25, 4).Value = 0.1866898467640965
$ws.Cells.Item(25, 5).Value = 0.07220609419841395
$ws.Cells.Item(25, 6).Value = 5.833562558737299
$ws.Cells.Item(25, 7).Value = 0.002579196904146021
$ws.Cells.Item(25, 9).Value = 4.143116904000891
$ws.Cells.Item(25, 10).Value = 0.07566807101770223
$ws.Cells.Item(25, 13).Value = 3.481240957586976
$ws.Cells.Item(25, 14).Value = 1.544987647507867
